$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Job adverts by occupation" row (row 13) to reflect June 2025 data.
# Set D13 (Next period) before C13 (Latest period) so new shared strings are
# appended in the same order as the target workbook.
$ws.Range("D13").Value = "Jul 2025 (29/08/25)"
$ws.Range("C13").Value = "Jun 2025 (01/08/25)"

# Update the selected cell shown when the workbook is opened.
$ws.Range("C13").Select()
